# Re-window the reaction-sensitivity tables on both sheets: the
# "fixed workflow" now starts sampling 4 cutoff steps later, so the
# first 4 original data rows are dropped, the remaining rows shift up
# to the top, and the freed-up trailing rows are removed outright.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Snapshot the existing 19 data rows (rows 2..20) before anything
    # is overwritten, since values get copied upward by 4 rows.
    $origB = @{}
    $origC = @{}
    for ($r = 2; $r -le 20; $r++) {
        $origB[$r] = $ws.Cells.Item($r, 2).Value2
        $origC[$r] = $ws.Cells.Item($r, 3).Value2
    }

    # Column A (the Cutoff index, 0..14) is unchanged for the rows that
    # remain; columns B (Reaction step) and C (value) are refreshed
    # from the data that used to live 4 rows further down.
    for ($r = 2; $r -le 16; $r++) {
        $srcRow = $r + 4
        $ws.Cells.Item($r, 2).Value = $origB[$srcRow]
        $ws.Cells.Item($r, 3).Value = $origC[$srcRow]
    }

    # Remove the now-obsolete trailing rows 17..20 entirely (shrinks
    # the used range / dimension down to A1:C16).
    $ws.Range("A17:C20").EntireRow.Delete()
}
